$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: was sourced from old row 12
$ws.Range("D2").Value = 44358
$ws.Range("M2").Value = 100
$ws.Range("N2").Value = 11000
$ws.Range("O2").Value = 12000
$ws.Range("P2").Value = 11500
$ws.Range("Q2").Value = "$/caja 18 kilos granel"
$ws.Range("S2").Value = 639

# Row 3: was sourced from old row 4
$ws.Range("D3").Value = 44425
$ws.Range("L3").Value = "Primera"
$ws.Range("M3").Value = 100
$ws.Range("N3").Value = 12000
$ws.Range("O3").Value = 13000
$ws.Range("P3").Value = 12500
$ws.Range("S3").Value = 694

# Row 4: was sourced from old row 8
$ws.Range("D4").Value = 44272
$ws.Range("N4").Value = 9000
$ws.Range("O4").Value = 10000
$ws.Range("P4").Value = 9500
$ws.Range("Q4").Value = "$/caja 15 kilos granel"
$ws.Range("S4").Value = 633
$ws.Range("T4").Value = 15

# Row 5: was sourced from old row 9
$ws.Range("D5").Value = 44272
$ws.Range("L5").Value = "Segunda"
$ws.Range("N5").Value = 8000
$ws.Range("O5").Value = 8000
$ws.Range("P5").Value = 8000
$ws.Range("Q5").Value = "$/caja 15 kilos granel"
$ws.Range("S5").Value = 533
$ws.Range("T5").Value = 15

# Row 6: was sourced from old row 13
$ws.Range("D6").Value = 44776
$ws.Range("M6").Value = 50
$ws.Range("N6").Value = 10000
$ws.Range("P6").Value = 10000
$ws.Range("Q6").Value = "$/bandeja 18 kilos granel"
$ws.Range("S6").Value = 556
$ws.Range("T6").Value = 18

# Row 7: was sourced from old row 14
$ws.Range("D7").Value = 44776
$ws.Range("L7").Value = "Segunda"
$ws.Range("M7").Value = 50
$ws.Range("N7").Value = 8000
$ws.Range("O7").Value = 8000
$ws.Range("P7").Value = 8000
$ws.Range("Q7").Value = "$/bandeja 18 kilos granel"
$ws.Range("S7").Value = 444

# Row 8: was sourced from old row 5
$ws.Range("D8").Value = 44698
$ws.Range("M8").Value = 50
$ws.Range("N8").Value = 10000
$ws.Range("P8").Value = 10000
$ws.Range("Q8").Value = "$/caja 18 kilos granel"
$ws.Range("S8").Value = 556
$ws.Range("T8").Value = 18

# Row 9: was sourced from old row 10
$ws.Range("D9").Value = 44299
$ws.Range("L9").Value = "Primera"
$ws.Range("M9").Value = 100
$ws.Range("N9").Value = 10000
$ws.Range("O9").Value = 11000
$ws.Range("P9").Value = 10500
$ws.Range("Q9").Value = "$/caja 18 kilos granel"
$ws.Range("R9").Value = "Región del Maule"
$ws.Range("S9").Value = 583
$ws.Range("T9").Value = 18

# Row 10: was sourced from old row 11
$ws.Range("L10").Value = "Segunda"
$ws.Range("M10").Value = 50
$ws.Range("N10").Value = 9000
$ws.Range("O10").Value = 9000
$ws.Range("P10").Value = 9000
$ws.Range("S10").Value = 500

# Row 11: was sourced from old row 6
$ws.Range("D11").Value = 44363
$ws.Range("L11").Value = "Primera"
$ws.Range("M11").Value = 100
$ws.Range("O11").Value = 10000
$ws.Range("P11").Value = 9500
$ws.Range("Q11").Value = "$/caja 15 kilos empedrada"
$ws.Range("R11").Value = "Región de O'Higgins"
$ws.Range("S11").Value = 633
$ws.Range("T11").Value = 15

# Row 12: was sourced from old row 2
$ws.Range("D12").Value = 44307
$ws.Range("M12").Value = 50
$ws.Range("N12").Value = 10000
$ws.Range("O12").Value = 10000
$ws.Range("P12").Value = 10000
$ws.Range("Q12").Value = "$/bandeja 18 kilos granel"
$ws.Range("S12").Value = 556

# Row 13: was sourced from old row 3
$ws.Range("D13").Value = 44307
$ws.Range("L13").Value = "Segunda"
$ws.Range("N13").Value = 8000
$ws.Range("O13").Value = 8000
$ws.Range("P13").Value = 8000
$ws.Range("S13").Value = 444

# Row 14: was sourced from old row 7
$ws.Range("D14").Value = 44316
$ws.Range("L14").Value = "Primera"
$ws.Range("M14").Value = 100
$ws.Range("N14").Value = 9000
$ws.Range("O14").Value = 10000
$ws.Range("P14").Value = 9500
$ws.Range("Q14").Value = "$/caja 18 kilos granel"
$ws.Range("S14").Value = 528
